$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(4)

# Trim the remark down to a single line: drop the "warn()"/"info()" bullets
# and the ", also" qualifier, keeping just "Halts execution".
$sh.TextFrame.TextRange.Text = "Halts execution"

# Re-anchor/resize the box to its new (smaller, one-line) footprint.
$sh.Left = 480.0
$sh.Top = 124.33275590551182
